$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cells value and force it to be stored as plain text
# (matching the source data, which uses text such as "413.86" or
# "  +4.29%  " rather than numeric/percentage cell types). Using
# TextToColumns with a Text column format re-parses the just-entered
# value as text without leaving behind any new/altered cell styles.
function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.Value = $val
    $cell.TextToColumns(1,1,1,1,1,1,1,1,1,2)
}

Set-TextValue "D2" '60.656.45'
Set-TextValue "E2" '  +6.77%  '
Set-TextValue "D3" '3.361.76'
Set-TextValue "E3" '  +3.11%  '
Set-TextValue "E4" '  -0.05%  '
Set-TextValue "D5" '413.86'
Set-TextValue "E5" '  +4.29%  '
Set-TextValue "D6" '112.61'
Set-TextValue "E6" '  +2.35%  '
Set-TextValue "D7" '0.588'
Set-TextValue "E7" '  +4.97%  '
Set-TextValue "E8" '  +0.08%  '
Set-TextValue "D9" '0.639'
Set-TextValue "E9" '  +2.28%  '
Set-TextValue "D10" '40.09'
Set-TextValue "E10" '  +2.18%  '
Set-TextValue "D11" '0.0990'
Set-TextValue "E11" '  +2.37%  '
Set-TextValue "D13" '3.889.88'
Set-TextValue "E13" '  +2.76%  '
Set-TextValue "D14" '8.58'
Set-TextValue "E14" '  +5.43%  '
Set-TextValue "D15" '20.28'
Set-TextValue "E15" '  +6.28%  '
Set-TextValue "D16" '3.326.64'
Set-TextValue "E16" '  +1.84%  '
Set-TextValue "E17" '  +0.98%  '
Set-TextValue "D18" '60.394.34'
Set-TextValue "E18" '  +6.44%  '
Set-TextValue "D19" '10.88'
Set-TextValue "E19" '  +0.21%  '
Set-TextValue "D20" '3.40'
Set-TextValue "E20" '  +2.88%  '
Set-TextValue "E21" '  +5.19%  '
Set-TextValue "D22" '13.27'
Set-TextValue "E22" '  +2.42%  '
Set-TextValue "D23" '305.27'
Set-TextValue "E23" '  -1.53%  '
Set-TextValue "D24" '75.84'
Set-TextValue "E24" '  +0.83%  '
Set-TextValue "D25" '3.21'
Set-TextValue "E25" '  +1.50%  '
Set-TextValue "D26" '28.83'
Set-TextValue "E26" '  +2.30%  '
Set-TextValue "E27" '  +2.48%  '
Set-TextValue "D28" '0.181'
Set-TextValue "E28" '  +6.94%  '
Set-TextValue "D29" '7.99'
Set-TextValue "E29" '  +1.03%  '
Set-TextValue "D30" '7.54'
Set-TextValue "E30" '  +3.90%  '
Set-TextValue "B31" 'Hedera'
Set-TextValue "C31" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D31" '0.115'
Set-TextValue "E31" '  +5.20%  '
Set-TextValue "B32" 'Toncoin'
Set-TextValue "C32" 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D32" '2.64'
Set-TextValue "E32" '  +23.66%  '
Set-TextValue "D33" '11.63'
Set-TextValue "E33" '  +5.68%  '
Set-TextValue "D34" '0.999'
Set-TextValue "E34" '  +0.27%  '
Set-TextValue "D35" '40.18'
Set-TextValue "E35" '  +6.98%  '
Set-TextValue "E36" '  +6.30%  '
Set-TextValue "D37" '52.34'
Set-TextValue "E37" '  +1.65%  '
Set-TextValue "D38" '3.14'
Set-TextValue "E38" '  +1.32%  '
Set-TextValue "D39" '0.998'
Set-TextValue "E39" '  -0.31%  '
Set-TextValue "E40" '  -2.99%  '
Set-TextValue "D41" '138.02'
Set-TextValue "E41" '  +2.16%  '
Set-TextValue "E42" '  +2.54%  '
Set-TextValue "D43" '1.93'
Set-TextValue "E43" '  +0.31%  '
Set-TextValue "B44" 'NEARProtocol'
Set-TextValue "C44" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D44" '3.99'
Set-TextValue "E44" '  +0.45%  '
Set-TextValue "B45" 'TheGraph'
Set-TextValue "C45" 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue "D45" '0.290'
Set-TextValue "E45" '  +3.53%  '
Set-TextValue "D46" '17.01'
Set-TextValue "E46" '  -1.69%  '
Set-TextValue "D47" '2.29'
Set-TextValue "E47" '  +9.47%  '
Set-TextValue "D48" '22.51'
Set-TextValue "E48" '  +2.23%  '
Set-TextValue "D49" '2.209.50'
Set-TextValue "E49" '  +2.77%  '
Set-TextValue "D50" '2.41'
Set-TextValue "E50" '  +1.35%  '
Set-TextValue "D51" '1.98'
Set-TextValue "E51" '  -2.03%  '
